{"js": "// Update the date line and all two-digit x two-digit multiplication\n// answers in the table to the values from the latest generated output.\nconst replacements = [\n  [\"2024-11-22 Friday\", \"2024-11-23 Saturday\"],\n  [\"54\u00d748=2592\", \"22\u00d766=1452\"],\n  [\"22\u00d774=1628\", \"82\u00d798=8036\"],\n  [\"44\u00d794=4136\", \"38\u00d788=3344\"],\n  [\"96\u00d784=8064\", \"30\u00d741=1230\"],\n  [\"46\u00d729=1334\", \"73\u00d740=2920\"],\n  [\"19\u00d735=665\", \"16\u00d787=1392\"],\n  [\"41\u00d712=492\", \"98\u00d737=3626\"],\n  [\"66\u00d761=4026\", \"57\u00d798=5586\"],\n  [\"31\u00d756=1736\", \"59\u00d725=1475\"],\n  [\"81\u00d771=5751\", \"86\u00d782=7052\"],\n  [\"36\u00d793=3348\", \"94\u00d794=8836\"],\n  [\"57\u00d788=5016\", \"83\u00d774=6142\"],\n  [\"37\u00d731=1147\", \"11\u00d731=341\"],\n  [\"23\u00d791=2093\", \"12\u00d757=684\"],\n  [\"11\u00d797=1067\", \"71\u00d766=4686\"],\n  [\"33\u00d797=3201\", \"11\u00d746=506\"],\n  [\"68\u00d753=3604\", \"68\u00d719=1292\"],\n  [\"40\u00d767=2680\", \"97\u00d739=3783\"],\n  [\"97\u00d781=7857\", \"39\u00d756=2184\"],\n  [\"65\u00d754=3510\", \"13\u00d796=1248\"],\n  [\"41\u00d730=1230\", \"42\u00d780=3360\"],\n  [\"66\u00d754=3564\", \"87\u00d766=5742\"],\n  [\"21\u00d713=273\", \"67\u00d743=2881\"],\n  [\"82\u00d754=4428\", \"92\u00d740=3680\"],\n  [\"77\u00d797=7469\", \"74\u00d740=2960\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all two-digit x two-digit multiplication\n# answers in the table to the values from the latest generated output.\n$d = $word.ActiveDocument\n$pairs = @(\n    @(\"2024-11-22 Friday\", \"2024-11-23 Saturday\"),\n    @(\"54\u00d748=2592\", \"22\u00d766=1452\"),\n    @(\"22\u00d774=1628\", \"82\u00d798=8036\"),\n    @(\"44\u00d794=4136\", \"38\u00d788=3344\"),\n    @(\"96\u00d784=8064\", \"30\u00d741=1230\"),\n    @(\"46\u00d729=1334\", \"73\u00d740=2920\"),\n    @(\"19\u00d735=665\", \"16\u00d787=1392\"),\n    @(\"41\u00d712=492\", \"98\u00d737=3626\"),\n    @(\"66\u00d761=4026\", \"57\u00d798=5586\"),\n    @(\"31\u00d756=1736\", \"59\u00d725=1475\"),\n    @(\"81\u00d771=5751\", \"86\u00d782=7052\"),\n    @(\"36\u00d793=3348\", \"94\u00d794=8836\"),\n    @(\"57\u00d788=5016\", \"83\u00d774=6142\"),\n    @(\"37\u00d731=1147\", \"11\u00d731=341\"),\n    @(\"23\u00d791=2093\", \"12\u00d757=684\"),\n    @(\"11\u00d797=1067\", \"71\u00d766=4686\"),\n    @(\"33\u00d797=3201\", \"11\u00d746=506\"),\n    @(\"68\u00d753=3604\", \"68\u00d719=1292\"),\n    @(\"40\u00d767=2680\", \"97\u00d739=3783\"),\n    @(\"97\u00d781=7857\", \"39\u00d756=2184\"),\n    @(\"65\u00d754=3510\", \"13\u00d796=1248\"),\n    @(\"41\u00d730=1230\", \"42\u00d780=3360\"),\n    @(\"66\u00d754=3564\", \"87\u00d766=5742\"),\n    @(\"21\u00d713=273\", \"67\u00d743=2881\"),\n    @(\"82\u00d754=4428\", \"92\u00d740=3680\"),\n    @(\"77\u00d797=7469\", \"74\u00d740=2960\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: not found\"\n    }\n}\n"}
